# BOT; UPDATE DATA
# Appends the 2020-05-13 (serial 43964) daily report row to the "all",
# "kobe" and "other" sheets, pushing the trailing footnote row down by one,
# and corrects two previously-provisional figures on "kobe" row 90.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": new row 36, footnote moves from row 36 -> row 37
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(36).Insert(-4121) | Out-Null

$wsAll.Range("A36").Value = 43964
$wsAll.Range("B36").Value = 278
$wsAll.Range("C36").Value = 276
$wsAll.Range("D36").Value = 79
$wsAll.Range("E36").Value = 68
$wsAll.Range("F36").Value = 11
$wsAll.Range("G36").Value = 9
$wsAll.Range("H36").Value = 188

$wsAll.Range("B37").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "kobe": correct row 90 (D/E), new row 91, footnote moves 91 -> 92
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D90").Value = 1
$wsKobe.Range("E90").Value = 279

$wsKobe.Rows.Item(91).Insert(-4121) | Out-Null

$wsKobe.Range("A91").Value = 43964
$wsKobe.Range("C91").Value = 2713
$wsKobe.Range("D91").Value = 2
$wsKobe.Range("E91").Value = 281
$wsKobe.Range("F91").Value = 74
$wsKobe.Range("G91").Value = 64
$wsKobe.Range("H91").Value = 10
$wsKobe.Range("I91").Value = 9
$wsKobe.Range("J91").Value = 179

$wsKobe.Range("K91").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "other": new row 66, footnote moves from row 66 -> row 67
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(66).Insert(-4121) | Out-Null

$wsOther.Range("A66").Value = 43964
$wsOther.Range("B66").Value = 0
$wsOther.Range("C66").Value = 14
$wsOther.Range("D66").Value = 5
$wsOther.Range("E66").Value = 4
$wsOther.Range("F66").Value = 1
$wsOther.Range("G66").Value = 0
$wsOther.Range("H66").Value = 9

$wsOther.Range("B67").Select() | Out-Null

$wsAll.Activate() | Out-Null
